$d = $word.ActiveDocument

# Every header/footer part in this document contains a single paragraph
# whose pPr explicitly points at a built-in "Header"/"Footer" style via
# <w:pStyle w:val="Header"/> or <w:pStyle w:val="Footer"/>. That style is
# not actually defined in styles.xml, which corrupts the document. Repair
# this by clearing each header/footer paragraph's explicit style back to
# the document's default ("Normal"), which drops the stray <w:pStyle>
# reference from every header and footer part (default, first-page and
# even-page, in every section).
for ($i = 1; $i -le $d.Sections.Count; $i++) {
    $sec = $d.Sections.Item($i)

    for ($j = 1; $j -le $sec.Headers.Count; $j++) {
        $sec.Headers.Item($j).Range.Style = "Normal"
    }

    for ($j = 1; $j -le $sec.Footers.Count; $j++) {
        $sec.Footers.Item($j).Range.Style = "Normal"
    }
}
